$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the two blank placeholder rows (13:14), which shifts the "Jan to Mar 2024"
# block (old rows 15-17) up to rows 13-15.
$ws.Rows("13:14").Delete()

# Delete one more blank row from the gap that separated the last quarterly block
# from the grand-total block, shifting the totals block (now rows 19-20) up to rows 18-19.
$ws.Rows("16").Delete()

# Column G picks up the same width/formatting as column F.
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# Move/refresh the active selection to A20, matching where the cursor ended up.
$ws.Range("A20").Select()
